$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the streetView column (Z) cell text and hyperlink targets with the
# refreshed Google Maps Street View embed links.

$r = $ws.Range("Z2")
$r.Hyperlinks.Delete()
$r.Value2 = 'https://www.google.com/maps/embed?pb=!4v1592665065170!6m8!1m7!1sMouHbQNkfq1KxzhNbKlxfg!2m2!1d36.89904235162321!2d-76.3114980328535!3f157.09!4f10.400000000000006!5f0.8160813932612223'

$r = $ws.Range("Z3")
$r.Value2 = 'https://www.google.com/maps/embed?pb=!4v1592665172692!6m8!1m7!1s-tCUM7tirQ1gU2se5Vkymw!2m2!1d36.87941427509518!2d-76.30873123639898!3f171.66!4f10!5f0.8160813932612223'
$r.Hyperlinks.Item(1).Address = 'https://www.google.com/maps/embed?pb=!4v1592665172692!6m8!1m7!1s-tCUM7tirQ1gU2se5Vkymw!2m2!1d36.87941427509518!2d-76.30873123639898!3f171.66!4f10!5f0.8160813932612223'

$r = $ws.Range("Z10")
$r.Value2 = 'https://www.google.com/maps/embed?pb=!4v1592665234801!6m8!1m7!1sR00yvp93G_GbsISBr5Uxqg!2m2!1d36.87457631912579!2d-76.29724063682792!3f289.46395626850057!4f-6.911394967165293!5f0.4000000000000002'
$r.Hyperlinks.Item(1).Address = 'https://www.google.com/maps/embed?pb=!4v1592665234801!6m8!1m7!1sR00yvp93G_GbsISBr5Uxqg!2m2!1d36.87457631912579!2d-76.29724063682792!3f289.46395626850057!4f-6.911394967165293!5f0.4000000000000002'

$r = $ws.Range("Z14")
$r.Value2 = 'https://www.google.com/maps/embed?pb=!4v1592665335209!6m8!1m7!1sO5Egiw4IOB1_r1tu4BWFwA!2m2!1d36.85804139313613!2d-76.27273965538782!3f116.88!4f10.120000000000005!5f0.4000000000000002'
$r.Hyperlinks.Item(1).Address = 'https://www.google.com/maps/embed?pb=!4v1592665335209!6m8!1m7!1sO5Egiw4IOB1_r1tu4BWFwA!2m2!1d36.85804139313613!2d-76.27273965538782!3f116.88!4f10.120000000000005!5f0.4000000000000002'

$r = $ws.Range("Z15")
$r.Value2 = 'https://www.google.com/maps/embed?pb=!4v1592665502514!6m8!1m7!1seBOvGINkRVe6370BDQEVbA!2m2!1d36.86314790329126!2d-76.27614736218922!3f168.34985944483896!4f5.958496962112932!5f0.772798507860903'
$r.Hyperlinks.Item(1).Address = 'https://www.google.com/maps/embed?pb=!4v1592665502514!6m8!1m7!1seBOvGINkRVe6370BDQEVbA!2m2!1d36.86314790329126!2d-76.27614736218922!3f168.34985944483896!4f5.958496962112932!5f0.772798507860903'

$r = $ws.Range("Z16")
$r.Value2 = 'https://www.google.com/maps/embed?pb=!4v1592665674872!6m8!1m7!1sCIA3gXXlK8Erlc2ceG4Vhw!2m2!1d36.86350246537197!2d-76.27355274908923!3f165.71932234716053!4f9.303222430462228!5f0.4000000000000002'
$r.Hyperlinks.Item(1).Address = 'https://www.google.com/maps/embed?pb=!4v1592665674872!6m8!1m7!1sCIA3gXXlK8Erlc2ceG4Vhw!2m2!1d36.86350246537197!2d-76.27355274908923!3f165.71932234716053!4f9.303222430462228!5f0.4000000000000002'

$r = $ws.Range("Z17")
$r.Value2 = 'https://www.google.com/maps/embed?pb=!4v1592665715906!6m8!1m7!1swxl7W_svb4ZWQ9m0EalZxw!2m2!1d36.86348702612092!2d-76.27366095485304!3f163.41!4f10!5f0.8160813932612223'
$r.Hyperlinks.Item(1).Address = 'https://www.google.com/maps/embed?pb=!4v1592665715906!6m8!1m7!1swxl7W_svb4ZWQ9m0EalZxw!2m2!1d36.86348702612092!2d-76.27366095485304!3f163.41!4f10!5f0.8160813932612223'

$r = $ws.Range("Z18")
$r.Value2 = 'https://www.google.com/maps/embed?pb=!4v1592665802086!6m8!1m7!1siBiQv1U88xQo56FPN0bTbw!2m2!1d36.8653228629351!2d-76.27600445696524!3f354.71864308642256!4f17.18177817873854!5f0.6222866636070421'
$r.Hyperlinks.Item(1).Address = 'https://www.google.com/maps/embed?pb=!4v1592665802086!6m8!1m7!1siBiQv1U88xQo56FPN0bTbw!2m2!1d36.8653228629351!2d-76.27600445696524!3f354.71864308642256!4f17.18177817873854!5f0.6222866636070421'

$r = $ws.Range("Z19")
$r.Value2 = 'https://www.google.com/maps/embed?pb=!4v1592665866558!6m8!1m7!1sK8CIvENU77mgE5vtfrZV1Q!2m2!1d36.86407357873125!2d-76.27335722979419!3f226.3629951882023!4f6.272487650330447!5f0.40457075809216086'
$r.Hyperlinks.Item(1).Address = 'https://www.google.com/maps/embed?pb=!4v1592665866558!6m8!1m7!1sK8CIvENU77mgE5vtfrZV1Q!2m2!1d36.86407357873125!2d-76.27335722979419!3f226.3629951882023!4f6.272487650330447!5f0.40457075809216086'

$r = $ws.Range("Z20")
$r.Value2 = 'https://www.google.com/maps/embed?pb=!4v1592665910517!6m8!1m7!1sEOjj2zmBL_8aVmC6UzbDXA!2m2!1d36.86547768723604!2d-76.27491296817067!3f148.13217185869334!4f3.43206204996757!5f0.8160813932612223'
$r.Hyperlinks.Item(1).Address = 'https://www.google.com/maps/embed?pb=!4v1592665910517!6m8!1m7!1sEOjj2zmBL_8aVmC6UzbDXA!2m2!1d36.86547768723604!2d-76.27491296817067!3f148.13217185869334!4f3.43206204996757!5f0.8160813932612223'

$r = $ws.Range("Z21")
$r.Value2 = 'https://www.google.com/maps/embed?pb=!4v1592665965000!6m8!1m7!1sQ1I0oMD2odkZEZ1PXIB6ZA!2m2!1d36.86452796795547!2d-76.2723941379832!3f45.99038225872033!4f12.969750919508996!5f0.8160813932612223'
$r.Hyperlinks.Item(1).Address = 'https://www.google.com/maps/embed?pb=!4v1592665965000!6m8!1m7!1sQ1I0oMD2odkZEZ1PXIB6ZA!2m2!1d36.86452796795547!2d-76.2723941379832!3f45.99038225872033!4f12.969750919508996!5f0.8160813932612223'

$r = $ws.Range("Z22")
$r.Value2 = 'https://www.google.com/maps/embed?pb=!4v1592666012188!6m8!1m7!1sfPG6cEeUFEdTGJomKGKTng!2m2!1d36.86376641053696!2d-76.27121155879614!3f41.27!4f10!5f0.8160813932612223'
$r.Hyperlinks.Item(1).Address = 'https://www.google.com/maps/embed?pb=!4v1592666012188!6m8!1m7!1sfPG6cEeUFEdTGJomKGKTng!2m2!1d36.86376641053696!2d-76.27121155879614!3f41.27!4f10!5f0.8160813932612223'

$r = $ws.Range("Z23")
$r.Value2 = 'https://www.google.com/maps/embed?pb=!4v1592666051032!6m8!1m7!1sCG3Ng0Hx8RZEQ13yLzrYig!2m2!1d36.85762740240267!2d-76.26594294015538!3f105.73!4f10!5f0.8160813932612223'
$r.Hyperlinks.Item(1).Address = 'https://www.google.com/maps/embed?pb=!4v1592666051032!6m8!1m7!1sCG3Ng0Hx8RZEQ13yLzrYig!2m2!1d36.85762740240267!2d-76.26594294015538!3f105.73!4f10!5f0.8160813932612223'

$r = $ws.Range("Z24")
$r.Value2 = 'https://www.google.com/maps/embed?pb=!4v1592666090791!6m8!1m7!1sT_es5p9pD31t-YdcIFDQ3w!2m2!1d36.85262526672637!2d-76.25288425539694!3f93.08!4f5.280000000000001!5f0.4000000000000002'
$r.Hyperlinks.Item(1).Address = 'https://www.google.com/maps/embed?pb=!4v1592666090791!6m8!1m7!1sT_es5p9pD31t-YdcIFDQ3w!2m2!1d36.85262526672637!2d-76.25288425539694!3f93.08!4f5.280000000000001!5f0.4000000000000002'

$r = $ws.Range("Z25")
$r.Value2 = 'https://www.google.com/maps/embed?pb=!4v1592666132074!6m8!1m7!1sRuz1hRa0iLnmbPUVE859-w!2m2!1d36.85231287239939!2d-76.25259655463324!3f175.0958478974717!4f13.479749444927066!5f0.4000000000000002'
$r.Hyperlinks.Item(1).Address = 'https://www.google.com/maps/embed?pb=!4v1592666132074!6m8!1m7!1sRuz1hRa0iLnmbPUVE859-w!2m2!1d36.85231287239939!2d-76.25259655463324!3f175.0958478974717!4f13.479749444927066!5f0.4000000000000002'

$r = $ws.Range("Z28")
$r.Value2 = 'https://www.google.com/maps/embed?pb=!4v1592666601823!6m8!1m7!1sCb2GYi55avnnrzEDyGWizQ!2m2!1d36.86315019098656!2d-76.25047304898814!3f80.03128482601167!4f16.304269243189196!5f0.8160813932612223'
$r.Hyperlinks.Item(1).Address = 'https://www.google.com/maps/embed?pb=!4v1592666601823!6m8!1m7!1sCb2GYi55avnnrzEDyGWizQ!2m2!1d36.86315019098656!2d-76.25047304898814!3f80.03128482601167!4f16.304269243189196!5f0.8160813932612223'

$r = $ws.Range("Z29")
$r.Value2 = 'https://www.google.com/maps/embed?pb=!4v1592666660733!6m8!1m7!1szX0TfVNGZLmCf7OzWxMAdg!2m2!1d36.84289222935679!2d-76.23597456947694!3f190.29873437114486!4f1.9053999445556684!5f0.8160813932612223'
$r.Hyperlinks.Item(1).Address = 'https://www.google.com/maps/embed?pb=!4v1592666660733!6m8!1m7!1szX0TfVNGZLmCf7OzWxMAdg!2m2!1d36.84289222935679!2d-76.23597456947694!3f190.29873437114486!4f1.9053999445556684!5f0.8160813932612223'

$r = $ws.Range("Z30")
$r.Value2 = 'https://www.google.com/maps/embed?pb=!4v1592666719772!6m8!1m7!1sgefwOVajIVpFF49p8o8dLw!2m2!1d36.87819215425838!2d-76.2480911662013!3f196.67710380758538!4f10.889840812576509!5f0.4000000000000002'
$r.Hyperlinks.Item(1).Address = 'https://www.google.com/maps/embed?pb=!4v1592666719772!6m8!1m7!1sgefwOVajIVpFF49p8o8dLw!2m2!1d36.87819215425838!2d-76.2480911662013!3f196.67710380758538!4f10.889840812576509!5f0.4000000000000002'

$r = $ws.Range("Z31")
$r.Value2 = 'https://www.google.com/maps/embed?pb=!4v1592666778575!6m8!1m7!1sGdcC0h8tp6KS7HmIkGOAQQ!2m2!1d36.91809792339702!2d-76.26426239307821!3f338.4674807008512!4f11.219220755293307!5f0.7876917578202689'
$r.Hyperlinks.Item(1).Address = 'https://www.google.com/maps/embed?pb=!4v1592666778575!6m8!1m7!1sGdcC0h8tp6KS7HmIkGOAQQ!2m2!1d36.91809792339702!2d-76.26426239307821!3f338.4674807008512!4f11.219220755293307!5f0.7876917578202689'

$r = $ws.Range("Z32")
$r.Value2 = 'https://www.google.com/maps/embed?pb=!4v1592666820888!6m8!1m7!1syOXh5P2YAc0EgmMU_szb3Q!2m2!1d36.91809581220794!2d-76.26349290489136!3f179.9333766983505!4f9.991950961987811!5f0.5586361636509503'
$r.Hyperlinks.Item(1).Address = 'https://www.google.com/maps/embed?pb=!4v1592666820888!6m8!1m7!1syOXh5P2YAc0EgmMU_szb3Q!2m2!1d36.91809581220794!2d-76.26349290489136!3f179.9333766983505!4f9.991950961987811!5f0.5586361636509503'

$r = $ws.Range("Z33")
$r.Value2 = 'https://www.google.com/maps/embed?pb=!4v1592666877268!6m8!1m7!1sGG_LLEJcrlPic7Dv-GbKmw!2m2!1d36.91812480341552!2d-76.2630478571323!3f189.25068923950414!4f6.130099650796822!5f0.415847888206226'
$r.Hyperlinks.Item(1).Address = 'https://www.google.com/maps/embed?pb=!4v1592666877268!6m8!1m7!1sGG_LLEJcrlPic7Dv-GbKmw!2m2!1d36.91812480341552!2d-76.2630478571323!3f189.25068923950414!4f6.130099650796822!5f0.415847888206226'

$r = $ws.Range("Z34")
$r.Value2 = 'https://www.google.com/maps/embed?pb=!4v1592666912681!6m8!1m7!1sEEsgn4M4w19lv5nJF3RpeQ!2m2!1d36.91729288221892!2d-76.26423205359693!3f1.11!4f10!5f0.8160813932612223'
$r.Hyperlinks.Item(1).Address = 'https://www.google.com/maps/embed?pb=!4v1592666912681!6m8!1m7!1sEEsgn4M4w19lv5nJF3RpeQ!2m2!1d36.91729288221892!2d-76.26423205359693!3f1.11!4f10!5f0.8160813932612223'

$r = $ws.Range("Z35")
$r.Value2 = 'https://www.google.com/maps/embed?pb=!4v1592666950210!6m8!1m7!1skJY0XlTFcvK4V4NYPbrEYA!2m2!1d36.91725404602425!2d-76.26384721333136!3f6.87!4f19.58!5f0.4000000000000002'
$r.Hyperlinks.Item(1).Address = 'https://www.google.com/maps/embed?pb=!4v1592666950210!6m8!1m7!1skJY0XlTFcvK4V4NYPbrEYA!2m2!1d36.91725404602425!2d-76.26384721333136!3f6.87!4f19.58!5f0.4000000000000002'

$r = $ws.Range("Z36")
$r.Value2 = 'https://www.google.com/maps/embed?pb=!4v1592667188917!6m8!1m7!1sJMPFSuRfSSvseozxtvMYtw!2m2!1d36.91714325084339!2d-76.26257363155776!3f194.96679944940217!4f10.082113883271404!5f0.4516936881563587'
$r.Hyperlinks.Item(1).Address = 'https://www.google.com/maps/embed?pb=!4v1592667188917!6m8!1m7!1sJMPFSuRfSSvseozxtvMYtw!2m2!1d36.91714325084339!2d-76.26257363155776!3f194.96679944940217!4f10.082113883271404!5f0.4516936881563587'

$r = $ws.Range("Z38")
$r.Value2 = 'https://www.google.com/maps/embed?pb=!4v1592667227532!6m8!1m7!1s_N82mesAxnrX91JXeeRgeQ!2m2!1d36.91679721300318!2d-76.2951005920984!3f178.0699692599674!4f10.455795409990998!5f0.8160813932612223'
$r.Hyperlinks.Item(1).Address = 'https://www.google.com/maps/embed?pb=!4v1592667227532!6m8!1m7!1s_N82mesAxnrX91JXeeRgeQ!2m2!1d36.91679721300318!2d-76.2951005920984!3f178.0699692599674!4f10.455795409990998!5f0.8160813932612223'

$r = $ws.Range("Z40")
$r.Value2 = 'https://www.google.com/maps/embed?pb=!4v1592667268960!6m8!1m7!1s-ejLgWnhuKiKXgscOvPtPQ!2m2!1d36.8743018445842!2d-76.29140445175926!3f103.43!4f11.790000000000006!5f0.4000000000000002'
$r.Hyperlinks.Item(1).Address = 'https://www.google.com/maps/embed?pb=!4v1592667268960!6m8!1m7!1s-ejLgWnhuKiKXgscOvPtPQ!2m2!1d36.8743018445842!2d-76.29140445175926!3f103.43!4f11.790000000000006!5f0.4000000000000002'

$r = $ws.Range("Z43")
$r.Value2 = 'https://www.google.com/maps/embed?pb=!4v1592667311839!6m8!1m7!1scQWWfGHswCQSKBY0YVFrpg!2m2!1d36.87815689266255!2d-76.24798190046953!3f209.43!4f17.019999999999996!5f0.4000000000000002'
$r.Hyperlinks.Item(1).Address = 'https://www.google.com/maps/embed?pb=!4v1592667311839!6m8!1m7!1scQWWfGHswCQSKBY0YVFrpg!2m2!1d36.87815689266255!2d-76.24798190046953!3f209.43!4f17.019999999999996!5f0.4000000000000002'

$r = $ws.Range("Z46")
$r.Value2 = 'https://www.google.com/maps/embed?pb=!4v1592667376459!6m8!1m7!1s73rkJmwZ2ho91UnTOCd-KQ!2m2!1d36.91724356474509!2d-76.26373672751919!3f2.297877136421069!4f5.549315654572695!5f0.40435055844788387'
$r.Hyperlinks.Item(1).Address = 'https://www.google.com/maps/embed?pb=!4v1592667376459!6m8!1m7!1s73rkJmwZ2ho91UnTOCd-KQ!2m2!1d36.91724356474509!2d-76.26373672751919!3f2.297877136421069!4f5.549315654572695!5f0.40435055844788387'

$r = $ws.Range("Z49")
$r.Value2 = 'https://www.google.com/maps/embed?pb=!4v1592667410618!6m8!1m7!1swKtNmFyou0we7Rx09QZyYA!2m2!1d36.91617068400947!2d-76.26509212568918!3f38.83!4f10!5f0.8160813932612223'
$r.Hyperlinks.Item(1).Address = 'https://www.google.com/maps/embed?pb=!4v1592667410618!6m8!1m7!1swKtNmFyou0we7Rx09QZyYA!2m2!1d36.91617068400947!2d-76.26509212568918!3f38.83!4f10!5f0.8160813932612223'

$r = $ws.Range("Z50")
$r.Value2 = 'https://www.google.com/maps/embed?pb=!4v1592667450549!6m8!1m7!1sbsjZO5OJ8v4P0ZLRK-dbEg!2m2!1d36.91598266032241!2d-76.26413338658887!3f4.4!4f10!5f0.8160813932612223'
$r.Hyperlinks.Item(1).Address = 'https://www.google.com/maps/embed?pb=!4v1592667450549!6m8!1m7!1sbsjZO5OJ8v4P0ZLRK-dbEg!2m2!1d36.91598266032241!2d-76.26413338658887!3f4.4!4f10!5f0.8160813932612223'

$r = $ws.Range("Z51")
$r.Value2 = 'https://www.google.com/maps/embed?pb=!4v1592667492286!6m8!1m7!1sZLzWDRJdX5fCrRxvPna0AA!2m2!1d36.91602678859328!2d-76.26446360364872!3f9.73!4f10!5f0.8160813932612223'
$r.Hyperlinks.Item(1).Address = 'https://www.google.com/maps/embed?pb=!4v1592667492286!6m8!1m7!1sZLzWDRJdX5fCrRxvPna0AA!2m2!1d36.91602678859328!2d-76.26446360364872!3f9.73!4f10!5f0.8160813932612223'

# Restore the previously-active selection/view (no pinned top-left cell,
# selection on F19) on Sheet1.
$ws.Range("F19").Select()

